$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''51.956.23'
$ws.Range('E2').Value = '  +0.87%  '
$ws.Range('D3').Value = '''2.823.15'
$ws.Range('E3').Value = '  +3.00%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = '''355.38'
$ws.Range('E5').Value = '  +7.05%  '
$ws.Range('D6').Value = '''113.85'
$ws.Range('E6').Value = '  -1.67%  '
$ws.Range('D7').Value = '''0.552'
$ws.Range('E7').Value = '  +2.77%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('E9').Value = '  +5.11%  '
$ws.Range('D10').Value = '''42.02'
$ws.Range('E10').Value = '  +1.27%  '
$ws.Range('E11').Value = '  -0.23%  '
$ws.Range('D12').Value = '''20.07'
$ws.Range('E12').Value = '  -0.09%  '
$ws.Range('E13').Value = '  +1.23%  '
$ws.Range('D14').Value = '''7.71'
$ws.Range('E14').Value = '  +1.48%  '
$ws.Range('D15').Value = '''3.246.15'
$ws.Range('E15').Value = '  +2.66%  '
$ws.Range('D16').Value = '''2.837.74'
$ws.Range('E16').Value = '  +2.96%  '
$ws.Range('E17').Value = '  +1.92%  '
$ws.Range('D18').Value = '''51.888.44'
$ws.Range('E18').Value = '  +0.90%  '
$ws.Range('D19').Value = '''7.37'
$ws.Range('E19').Value = '  +7.92%  '
$ws.Range('D20').Value = '''3.16'
$ws.Range('E20').Value = '  -1.50%  '
$ws.Range('D21').Value = '''13.54'
$ws.Range('E21').Value = '  +1.11%  '
$ws.Range('D22').Value = '''0.0₃0998'
$ws.Range('E22').Value = '  +2.57%  '
$ws.Range('D23').Value = '''270.09'
$ws.Range('E23').Value = '  -2.90%  '
$ws.Range('D24').Value = '''69.71'
$ws.Range('D25').Value = '''2.80'
$ws.Range('E25').Value = '  +6.04%  '
$ws.Range('D26').Value = '''26.82'
$ws.Range('E26').Value = '  +0.30%  '
$ws.Range('E27').Value = '  +0.15%  '
$ws.Range('D30').Value = '''0.140'
$ws.Range('E30').Value = '  -0.25%  '
$ws.Range('D31').Value = '''50.81'
$ws.Range('E31').Value = '  +1.59%  '
$ws.Range('D32').Value = '''33.94'
$ws.Range('E32').Value = '  -3.02%  '
$ws.Range('D33').Value = '''0.0452'
$ws.Range('E33').Value = '  +31.75%  '
$ws.Range('D34').Value = '''5.84'
$ws.Range('E34').Value = '  +5.65%  '
$ws.Range('D35').Value = '''0.0831'
$ws.Range('E35').Value = '  +1.57%  '
$ws.Range('E36').Value = '  -0.09%  '
$ws.Range('E37').Value = '  +0.53%  '
$ws.Range('D38').Value = '''3.22'
$ws.Range('E38').Value = '  +0.58%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D39').Value = '''4.88'
$ws.Range('E39').Value = '  -2.32%  '
$ws.Range('B40').Value = 'Celestia'
$ws.Range('C40').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D40').Value = '''18.42'
$ws.Range('E40').Value = '  -3.38%  '
$ws.Range('B41').Value = 'EnergySwap'
$ws.Range('C41').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D41').Value = '''23.89'
$ws.Range('E41').Value = '  +3.97%  '
$ws.Range('B42').Value = 'Monero'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D42').Value = '''128.76'
$ws.Range('E42').Value = '  +0.89%  '
$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D43').Value = '''2.58'
$ws.Range('E43').Value = '  +6.80%  '
$ws.Range('E44').Value = '  +1.45%  '
$ws.Range('E45').Value = '  +0.62%  '
$ws.Range('D46').Value = '''3.36'
$ws.Range('E46').Value = '  +1.36%  '
$ws.Range('D47').Value = '''2.074.20'
$ws.Range('E47').Value = '  -0.58%  '
$ws.Range('E48').Value = '  +4.33%  '
$ws.Range('D49').Value = '''0.961'
$ws.Range('E49').Value = '  +9.97%  '
$ws.Range('D50').Value = '''5.69'
$ws.Range('E50').Value = '  +3.28%  '
$ws.Range('D51').Value = '''60.47'
$ws.Range('E51').Value = '  +1.00%  '
